# Update the 2D training schedule data (g-casa children task), no break screen.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (trialTrain, x_fixStart, y_fixStart, x_corrSteps, y_corrSteps,
#                x_nrSteps, y_nrSteps, alienID, praclen, version)
$data = @(
    @(1, 5, 3, 3, 7, -2, 4, 45, 5),
    @(2, 6, 2, 5, 7, -1, 5, 56, 5),
    @(3, 6, 4, 1, 5, -5, 1, 12, 5),
    @(4, 5, 1, 2, 4, -3, 3, 34, 5),
    @(5, 8, 4, 4, 6, -4, 2, 23, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

[void]$ws.Range("I1").Select()
